$wb = $excel.ActiveWorkbook

# --- LOGIN sheet: move selection from A3:B3 to B2 (drops tabSelected there) ---
$loginSheet = $wb.Worksheets.Item(1)
$loginSheet.Range("B2").Select() | Out-Null

# --- Add the new "3RD PARTY PROVIDER" sheet right after LOGIN ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$newSheet.Name = "3RD PARTY PROVIDER"

# --- Fill in the header row and the single data row, in the same order the
#     strings were originally authored (keeps sharedStrings ordering stable) ---
$newSheet.Range("A1").Value = "Provider Name"
$newSheet.Range("B1").Value = "Mobile Number"
$newSheet.Range("C1").Value = "Other Contact Number"
$newSheet.Range("D1").Value = "Email Address"
$newSheet.Range("A2").Value = "Saran "
$newSheet.Range("E1").Value = "Fax Number"
$newSheet.Range("F1").Value = "Vat Number"
$newSheet.Range("F2").Value = "27832GAF*&"
$newSheet.Range("E2").Value = "87378HS&$#"
$newSheet.Range("G1").Value = "Status"
$newSheet.Range("G2").Value = "Inactive"
$newSheet.Range("H1").Value = "Street"
$newSheet.Range("I1").Value = "Street2"
$newSheet.Range("J1").Value = "City"
$newSheet.Range("K1").Value = "Province"
$newSheet.Range("L1").Value = "Country"
$newSheet.Range("M1").Value = "Postal Code"
$newSheet.Range("N1").Value = "Comments"
$newSheet.Range("B2").Value = "97483HSF%#"
$newSheet.Range("C2").Value = "8727GTH&^%"
$newSheet.Range("D2").Value = "saran@gmail.com"

# --- Hyperlink the email address; Excel auto-applies the Hyperlink cell
#     style when a link is added, but the source file keeps D2 unstyled, so
#     reset it back to Normal afterwards ---
$newSheet.Hyperlinks.Add($newSheet.Range("D2"), "mailto:saran@gmail.com")
$newSheet.Range("D2").Style = "Normal"

# --- Column sizing to match the authored (bestFit) widths as closely as
#     the engine's width quantisation allows ---
$newSheet.Columns("A").ColumnWidth = 12.5
$newSheet.Columns("B").ColumnWidth = 13.166666666666666
$newSheet.Columns("C").ColumnWidth = 19.0
$newSheet.Columns("D").ColumnWidth = 13.666666666666666
$newSheet.Columns("E").ColumnWidth = 13.666666666666666
$newSheet.Columns("F").ColumnWidth = 10.0
$newSheet.Columns("M").ColumnWidth = 9.666666666666666
$newSheet.Columns("N").ColumnWidth = 9.0

# --- Page setup / selection on the new sheet ---
$newSheet.PageSetup.Orientation = 1
$newSheet.Range("B7").Select() | Out-Null
